# Refresh the Price (D), Volume(1h) (E), and Hora (G) columns for each symbol row
# with the latest scrape snapshot. Values are kept as literal text (NumberFormat "@")
# so they round-trip the same way the source feed writes them (e.g. "288.58", "-0.91%").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ 'D' = '288.58'; 'E' = '-0.91%'; 'G' = '10' }
    3 = @{ 'D' = '31.09'; 'E' = '1.09%'; 'G' = '10' }
    4 = @{ 'D' = '4.923'; 'E' = '-0.55%'; 'G' = '10' }
    5 = @{ 'D' = '0.07328'; 'E' = '1.52%'; 'G' = '10' }
    6 = @{ 'D' = '2.221'; 'E' = '20.70%'; 'G' = '10' }
    7 = @{ 'D' = '7.735'; 'E' = '0.59%'; 'G' = '10' }
    8 = @{ 'E' = '-1.01%'; 'G' = '10' }
    9 = @{ 'D' = '0.9024'; 'E' = '0.52%'; 'G' = '10' }
    10 = @{ 'D' = '0.09157'; 'E' = '19.21%'; 'G' = '10' }
    11 = @{ 'D' = '0.1698'; 'E' = '2.70%'; 'G' = '10' }
    12 = @{ 'D' = '0.08176'; 'E' = '2.85%'; 'G' = '10' }
    13 = @{ 'D' = '0.03124'; 'E' = '3.03%'; 'G' = '10' }
    14 = @{ 'D' = '0.09944'; 'E' = '-0.65%'; 'G' = '10' }
    15 = @{ 'E' = '-0.10%'; 'G' = '10' }
    16 = @{ 'D' = '0.005727'; 'E' = '0.82%'; 'G' = '10' }
    17 = @{ 'D' = '3.533'; 'E' = '1.80%'; 'G' = '10' }
    18 = @{ 'D' = '2.082'; 'E' = '0.07%'; 'G' = '10' }
    19 = @{ 'D' = '0.3329'; 'E' = '0.33%'; 'G' = '10' }
    20 = @{ 'E' = '0.01%'; 'G' = '10' }
    21 = @{ 'D' = '4.165'; 'E' = '3.01%'; 'G' = '10' }
    22 = @{ 'D' = '0.2101'; 'E' = '-11.98%'; 'G' = '10' }
    23 = @{ 'D' = '0.04527'; 'E' = '0.79%'; 'G' = '10' }
    24 = @{ 'E' = '-0.41%'; 'G' = '10' }
    25 = @{ 'D' = '0.004159'; 'E' = '3.78%'; 'G' = '10' }
    26 = @{ 'D' = '0.0001301'; 'E' = '4.00%'; 'G' = '10' }
    27 = @{ 'G' = '10' }
    28 = @{ 'G' = '10' }
    29 = @{ 'G' = '10' }
    30 = @{ 'G' = '10' }
    31 = @{ 'G' = '10' }
    32 = @{ 'G' = '10' }
    33 = @{ 'G' = '10' }
    34 = @{ 'G' = '10' }
    35 = @{ 'G' = '10' }
    36 = @{ 'G' = '10' }
    37 = @{ 'G' = '10' }
    38 = @{ 'G' = '10' }
    39 = @{ 'D' = '0.01575'; 'E' = '-1.01%'; 'G' = '10' }
    40 = @{ 'D' = '0.04454'; 'E' = '1.11%'; 'G' = '10' }
    41 = @{ 'D' = '0.007321'; 'E' = '0.65%'; 'G' = '10' }
    42 = @{ 'D' = '0.009540'; 'E' = '-5.38%'; 'G' = '10' }
    43 = @{ 'D' = '0.1328'; 'E' = '1.58%'; 'G' = '10' }
    44 = @{ 'D' = '0.002291'; 'E' = '13.76%'; 'G' = '10' }
    45 = @{ 'D' = '0.008323'; 'E' = '-12.51%'; 'G' = '10' }
    46 = @{ 'D' = '0.00006106'; 'E' = '2.08%'; 'G' = '10' }
    47 = @{ 'E' = '0.05%'; 'G' = '10' }
    48 = @{ 'D' = '2.242'; 'E' = '-0.23%'; 'G' = '10' }
    49 = @{ 'E' = '-33.28%'; 'G' = '10' }
    50 = @{ 'D' = '0.00002101'; 'E' = '0.05%'; 'G' = '10' }
    51 = @{ 'D' = '0.0002001'; 'E' = '0.05%'; 'G' = '10' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$row][$col]
    }
}
